# Commit: "adding full run for ZEV Jan R2-4 and modifying files for consistency in R2-4"
#
# The roboticS1Prep column (I) previously stored the text string "No" for
# every data row. For consistency, those values are converted to a real
# boolean FALSE value (so the column now holds native booleans instead of
# the text "No"). Once no cell references the "No" shared string anymore,
# Excel drops it from the shared-strings table, which is why every other
# shared-string index used later in the table shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows 2-41 all have roboticS1Prep ("No") in column I -> change to boolean FALSE
$ws.Range("I2:I41").Value = $false

# Update the last active selection recorded for the sheet
$null = $ws.Range("N20").Select()
